$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "L1cam", "Itga5", "ECs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 368.4011609839244, 3315.61044885532, 0.2147678609052866, 0.2147678609052865),
    @("ECs", "L1cam", "Itga5", "FAPs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 355.5827478194456, 3200.244730375011, 0.207295074532457, 0.2072950745324569),
    @("ECs", "L1cam", "Itga5", "M2", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 227.641022376438, 2048.769201387942, 0.1327085270293527, 0.1327085270293527),
    @("ECs", "L1cam", "Itga5", "sCs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 71.96127477068933, 647.651472936204, 0.04195146673599356, 0.04195146673599356),
    @("FAPs", "L1cam", "Itga5", "ECs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 11.73310456132889, 105.59794105196, 0.006840080964143934, 0.006840080964143932),
    @("FAPs", "L1cam", "Itga5", "FAPs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 11.32485453962033, 101.923690856583, 0.006602082300831607, 0.006602082300831606),
    @("FAPs", "L1cam", "Itga5", "M2", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 7.250074649214, 65.250671842926, 0.00422659640826507, 0.004226596408265068),
    @("FAPs", "L1cam", "Itga5", "sCs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 2.291874322534667, 20.626868902812, 0.001336100419444666, 0.001336100419444665),
    @("M2", "L1cam", "Itga5", "ECs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 127.7845581127733, 1150.06102301496, 0.07449492322258223, 0.07449492322258223),
    @("M2", "L1cam", "Itga5", "FAPs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 123.338330914262, 1110.044978228358, 0.07190289364815629, 0.07190289364815629),
    @("M2", "L1cam", "Itga5", "M2", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 78.96014055716398, 710.641265014476, 0.04603161520705103, 0.04603161520705103),
    @("M2", "L1cam", "Itga5", "sCs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 24.960669703768, 224.646027333912, 0.0145513918162586, 0.0145513918162586),
    @("sCs", "L1cam", "Itga5", "ECs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 109.4550746031244, 985.09567142812, 0.06380933267136846, 0.06380933267136843),
    @("sCs", "L1cam", "Itga5", "FAPs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 105.6466165476057, 950.819548928451, 0.06158910516788627, 0.06158910516788626),
    @("sCs", "L1cam", "Itga5", "M2", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 67.634056907958, 608.706512171622, 0.03942881636874755, 0.03942881636874754),
    @("sCs", "L1cam", "Itga5", "sCs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 21.38029825292934, 192.422684276364, 0.0124641326021745, 0.0124641326021745)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
